$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Thesis title: Global temperature response to volcanic activity" -> add
#    a period at the end, and insert a new paragraph right after it with the
#    PhD project description (kept as separate runs, mirroring the rest of
#    the document's run layout).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Thesis title: Global temperature response to volcanic activity",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Thesis title: Global temperature response to volcanic activity.", 2) | Out-Null

$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $candidate = $paras.Item($i)
    if ($candidate.Range.Text -like "Thesis title: Global temperature response to volcanic activity.*") {
        $thesisPara = $candidate
        break
    }
}

$afterRng = $thesisPara.Range
$afterRng.Collapse(0)
$afterRng.InsertParagraphAfter()

# Re-fetch the freshly created (still empty) paragraph that follows.
$newParaIndex = $i + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRng = $newPara.Range

$phdXml = '<?xml version="1.0" standalone="yes"?>' +
'<?mso-application progid="Word.Document"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p><w:pPr><w:pStyle w:val="Definition"/></w:pPr>' +
'<w:r><w:t xml:space="preserve">The PhD work consist of running long climate model simulations with volcanic forcing</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">and investigate the corresponding temperature response to volcanoes. The response to</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">volcanic forcing is hypothesized to be linear. Further, analysis will be carried out</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">to investigate the universality of the response to volcanic forcing with respect to</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">any kind of forcing, possibly providing valuble insight into the equilibrium climate</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">sensitivity.</w:t></w:r>' +
'</w:p>' +
'<w:sectPr/>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$newRng.InsertXML($phdXml)

# ---------------------------------------------------------------------------
# 2) Typo fix: "superthermal" -> "suprathermal"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "superthermal electrons observed by a moving radar numerically and compare to real",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "suprathermal electrons observed by a moving radar numerically and compare to real", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "within the terminal, for example with a file manager like" ->
#    "within the terminal, for example with the file manager"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "within the terminal, for example with a file manager like",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "within the terminal, for example with the file manager", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove " or nnn" (including the "nnn" hyperlink) so the sentence reads
#    "...with the file manager lf."
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.Address -eq "https://github.com/jarun/nnn") {
        $nnnLink = $candidate
        break
    }
}
$nnnRange = $nnnLink.Range
$removeRange = $d.Range($nnnRange.Start - 4, $nnnRange.End)
$removeRange.Delete()
